$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").BorderAround()
